$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - "How I Use Every Claude Code Feature"
$ws.Range("A6").Value = "How I Use Every Claude Code Feature"
$ws.Range("B6").Value = "A comprehensive practitioner's guide to every Claude Code feature, from CLAUDE.md configuration to hooks, planning mode, skills, MCP, and enterprise GitHub Actions workflows."
$ws.Range("C6").Value = "guide"
$ws.Range("D6").Value = "advanced"
# readingTime/createdAt columns store numeric/date-looking values as text
# (matching the rest of the sheet), so force text format before assignment.
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "15"
$ws.Range("F6").Value = "https://blog.sshh.io/p/how-i-use-every-claude-code-feature"
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = "skill"
$ws.Range("I6").Value = "All Features"
$ws.Range("J6").Value = "CLAUDE.md; hooks; planning mode; skills; MCP; custom commands; subagents; GitHub Actions; SDK; enterprise"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "2026-02-06"

# Row 7 - "Agentic Coding Recommendations — Armin Ronacher"
$ws.Range("A7").Value = "Agentic Coding Recommendations — Armin Ronacher"
$ws.Range("B7").Value = "Flask creator Armin Ronacher shares practical recommendations for agentic coding, covering language choices, tool design, observability, and code organization for AI-driven development."
$ws.Range("C7").Value = "guide"
$ws.Range("D7").Value = "intermediate"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "12"
$ws.Range("F7").Value = "https://lucumr.pocoo.org/2025/6/12/agentic-coding/"
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = "agent"
$ws.Range("I7").Value = "Agentic Coding"
$ws.Range("J7").Value = "agentic coding; Go; Python; tooling; observability; refactoring; workflow"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "2026-02-07"
